$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1 ("DataSet") - add a new "BillingDetails" row (row 17) plus two
# blank quote-prefixed cells on row 16 (L16/M16), mirroring the existing
# "AddressBook" row (14) layout.
# ---------------------------------------------------------------------------
$ws1.Range("L16").Value = "'"
$ws1.Range("L16").Value = ""
$ws1.Range("M16").Value = "'"
$ws1.Range("M16").Value = ""

$ws1.Range("A17").Value = "BillingDetails"
$ws1.Range("E17").Value = "Test"
$ws1.Range("F17").Value = "qa"
$ws1.Range("I17").Value = "6 Walnut Valley Dr"
$ws1.Range("J17").Value = "Little Rock"
$ws1.Range("K17").Value = "Arkansas"
$ws1.Range("L17").Value = "'72211"
$ws1.Range("M17").Value = "'9999999999"

# ---------------------------------------------------------------------------
# Sheet2 ("Sheet1") - add a new "Giftaccount" row (row 13).
# ---------------------------------------------------------------------------
$ws2.Range("B13").Value = "vickstest147@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("B13"), "mailto:vickstest147@gmail.com")
$ws2.Range("A13").Value = "Giftaccount"
$ws2.Range("C13").Value = "Ajitvv#1969"

$ws2.Range("H13").Value = "'"
$ws2.Range("H13").Value = ""

$ws2.Range("M13").Value = "'06492"

# ---------------------------------------------------------------------------
# Selections / view state to match the authored workbook.
# ---------------------------------------------------------------------------
$ws1.Range("B19").Select()
$ws2.Select()
$ws2.Range("J15").Select()
